$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 109; this shifts the existing rows 109-169 down
# to 110-170 and grows the sheet dimension from A1:T169 to A1:T170.
$ws.Rows("109").Insert()

# Populate the newly inserted row 109 with the new record.
$ws.Range("A109").Value = 10
$ws.Range("B109").Value = "Vega Modelo de Temuco"
$ws.Range("C109").Value = "La Araucanía"
$ws.Range("D109").Value = 44529
$ws.Range("E109").Value = 9
$ws.Range("F109").Value = "Fruta"
$ws.Range("G109").Value = 100102
$ws.Range("H109").Value = "Cítricos"
$ws.Range("I109").Value = 100102006
$ws.Range("J109").Value = "Pomelo"
$ws.Range("K109").Value = "Start Ruby"
$ws.Range("L109").Value = "Primera"
$ws.Range("M109").Value = 55
$ws.Range("N109").Value = 12000
$ws.Range("O109").Value = 12000
$ws.Range("P109").Value = 12000
$ws.Range("Q109").Value = "$/caja 14 kilos empedrada"
$ws.Range("R109").Value = "Región de O'Higgins"
$ws.Range("S109").Value = 857
$ws.Range("T109").Value = 14
